$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.679.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.47%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.642.86'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.82%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.93%  '

$ws.Range("E6").Value = '  +1.40%  '

$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  +0.78%  '

$ws.Range("E9").Value = '  +0.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.03'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.07%  '

$ws.Range("E11").Value = '  -0.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.871.86'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.648.00'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.81%  '

$ws.Range("E14").Value = '  +1.48%  '

$ws.Range("E15").Value = '  +1.45%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.98'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.705.46'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.01%  '

$ws.Range("E18").Value = '  +0.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.72'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.06%  '

$ws.Range("E20").Value = '  +0.00%  '

$ws.Range("E21").Value = '  +1.14%  '

$ws.Range("E22").Value = '  +1.24%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.49'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.58%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.25'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +14.40%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.48'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.92%  '

$ws.Range("E26").Value = '  -0.03%  '

$ws.Range("E27").Value = '  -0.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.16'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +4.42%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.69'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.98%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0516'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.80%  '

$ws.Range("E31").Value = '  +0.75%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.42%  '

$ws.Range("E33").Value = '  +2.76%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.278.38'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.94%  '

$ws.Range("E35").Value = '  +2.26%  '

$ws.Range("E36").Value = '  +1.31%  '

$ws.Range("E37").Value = '  +2.84%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.534'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.85%  '

$ws.Range("E39").Value = '  +3.88%  '

$ws.Range("E40").Value = '  +0.01%  '

$ws.Range("E41").Value = '  +2.34%  '

$ws.Range("E42").Value = '  -0.38%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.42%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.782.16'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.80%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '91.79'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.43%  '

$ws.Range("E47").Value = '  +1.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0516'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.93%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.76'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.87%  '

$ws.Range("E50").Value = '  +2.11%  '

$ws.Range("E51").Value = '  -0.35%  '
